$wb = $excel.ActiveWorkbook

# Mapping of worksheet name -> list of (cell address, expected old value, new value)
# This script updates column F ("想去人数" / interest count) values per the source diff.

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 107
$ws.Range("F4").Value = 8294
$ws.Range("F6").Value = 112
$ws.Range("F7").Value = 7251
$ws.Range("F8").Value = 1144
$ws.Range("F9").Value = 569
$ws.Range("F10").Value = 499
$ws.Range("F13").Value = 164
$ws.Range("F15").Value = 170
$ws.Range("F17").Value = 113
$ws.Range("F18").Value = 11941
$ws.Range("F21").Value = 150
$ws.Range("F22").Value = 2379
$ws.Range("F24").Value = 3378
$ws.Range("F25").Value = 52
$ws.Range("F27").Value = 2826
$ws.Range("F28").Value = 111
$ws.Range("F31").Value = 3314
$ws.Range("F33").Value = 2427
$ws.Range("F35").Value = 1676
$ws.Range("F37").Value = 114
$ws.Range("F38").Value = 5932
$ws.Range("F40").Value = 20
$ws.Range("F41").Value = 167
$ws.Range("F44").Value = 1108
$ws.Range("F45").Value = 1089
$ws.Range("F46").Value = 1563
$ws.Range("F48").Value = 109

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 109

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 288
$ws.Range("F3").Value = 430

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 107
$ws.Range("F4").Value = 288
$ws.Range("F5").Value = 430
$ws.Range("F8").Value = 8294
$ws.Range("F10").Value = 112
$ws.Range("F11").Value = 7251
$ws.Range("F12").Value = 7252
$ws.Range("F13").Value = 1144
$ws.Range("F14").Value = 569
$ws.Range("F15").Value = 499
$ws.Range("F18").Value = 164
$ws.Range("F20").Value = 113
$ws.Range("F22").Value = 11941
$ws.Range("F25").Value = 150
$ws.Range("F26").Value = 2379
$ws.Range("F27").Value = 2379
$ws.Range("F28").Value = 3378
$ws.Range("F29").Value = 2826
$ws.Range("F30").Value = 111
$ws.Range("F33").Value = 3314
$ws.Range("F36").Value = 2427
$ws.Range("F38").Value = 1676
$ws.Range("F39").Value = 114
$ws.Range("F40").Value = 5932
$ws.Range("F44").Value = 167
$ws.Range("F47").Value = 1108
$ws.Range("F48").Value = 1089
$ws.Range("F49").Value = 1563
$ws.Range("F50").Value = 109
